{"js": "// Find the run of text that ends the target paragraph and append a new\n// sentence right after it, inside the same paragraph, matching the\n// existing run's character formatting (b=0, bCs, sz=20, szCs=20).\nconst anchorText =\n  \"along with its corresponding tokenizer which will be included in the model architecture instead of done individually during the preprocessing stage.\";\n\nconst results = context.document.body.search(anchorText, {\n  matchCase: false,\n  matchWholeWord: false,\n});\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error(\"Anchor text not found: \" + anchorText);\n}\n\nconst hit = results.items[0];\nconst paragraph = hit.paragraphs.getFirst();\nconst tail = paragraph.getRange(\"End\");\n\ntail.insertText(\n  \"Along with the tokenizer, a sequence padder was added to fill in the empty tokens. \",\n  \"Before\"\n);\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# Locate the text (including its trailing space) that ends the \"Model\n# Architecture\" paragraph, then collapse to the end of that match and\n# insert the new sentence right there -- still inside the same paragraph,\n# after the existing text, inheriting its character formatting (not bold,\n# sz/szCs 20).\n$rng = $d.Content\n$find = $rng.Find\n$find.Text = \"preprocessing stage. \"\n$find.MatchCase = $true\n$find.MatchWholeWord = $false\n\nif ($find.Execute()) {\n    $rng.Collapse(0)  # wdCollapseEnd\n    $rng.InsertAfter(\"Along with the tokenizer, a sequence padder was added to fill in the empty tokens. \")\n}\n"}
